# v0.5.0 Historical data for aggregated months now importable from excel
#
# Adds a new "historical_scores" worksheet (month / score / day0 / ml)
# populated from the historical aggregated-month scoring data, and makes
# it the active sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Tidy up selection on the existing "monthly_targets" sheet before we
#    add/activate the new sheet (mirrors a user reviewing it, then moving
#    their selection back to the top-left data cell).
# ---------------------------------------------------------------------
$wsTargets = $wb.Worksheets.Item("monthly_targets")
$wsTargets.Activate()
$wsTargets.Range("A2").Select() | Out-Null

# ---------------------------------------------------------------------
# 2) Add the new sheet after the last existing sheet.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "historical_scores"

# ---------------------------------------------------------------------
# 3) Historical aggregated-month data: month, score, day0, ml.
#    Month values are written first (column A, top to bottom) and the
#    headers afterwards so new shared-string entries land in the same
#    order as the authored workbook.
# ---------------------------------------------------------------------
$data = @(
    @("18m07", 0.4767, 2, 626),
    @("18m08", 0.1422, 4, 573),
    @("18m09", 0.4659, 6, 450),
    @("18m10", 0.4454, 4, 590),
    @("18m11", 0.4806, 4, 546),
    @("18m12", 0.3931, 5, 608),
    @("19m01", 0.8286, 11, 382),
    @("19m02", 0.874, 10, 408),
    @("19m03", 0.6507, 5, 494),
    @("19m04", 1.1186, 8, 452),
    @("19m05", 1.121, 5, 511),
    @("19m06", 0.8249, 1, 627),
    @("19m07", 0.637, 3, 645),
    @("19m08", 0.7155, 3, 480),
    @("19m09", 1.015, 5, 402),
    @("19m10", 0.8672, 8, 472),
    @("19m11", 0.9352, 6, 475),
    @("19m12", 0.6356, 4, 519),
    @("20m01", 0.9185, 12, 314),
    @("20m02", 0.7231, 7, 402),
    @("20m03", 1.2712, 7, 408),
    @("20m04", 1.0843, 5, 542),
    @("20m05", 0.6307, 9, 397),
    @("20m06", 1.1801, 7, 456),
    @("20m07", 1.0217, 6, 434),
    @("20m08", 0.6323, 4, 599),
    @("20m09", 0.2471, 2, 513),
    @("20m10", 0.5243, 7, 397),
    @("20m11", 0.7656, 5, 475),
    @("20m12", 0.4579, 6, 513),
    @("21m01", 1.0493, 8, 421)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}

# Score column is a percentage.
$ws.Range("B2:B32").NumberFormat = "0.00%"

# Headers last, so "month"/"score"/"day0"/"ml" are appended to the shared
# string table after the month codes.
$ws.Range("A1").Value = "month"
$ws.Range("B1").Value = "score"
$ws.Range("C1").Value = "day0"
$ws.Range("D1").Value = "ml"

# ---------------------------------------------------------------------
# 4) View state for the new sheet: freeze the header row, zoom to 85%
#    and leave the selection on D16; this makes "historical_scores" the
#    active (selected) tab.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("D16").Select() | Out-Null
$excel.ActiveWindow.Zoom = 85
